# Pavel - new user for linking test
# Adds a new "Linking_AutoUser" row to the Users sheet for the Linking test suite.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 52), following the same layout as the existing rows:
# A=UserName  B=Password  C=ClientID  D=ROLE  E=Description  F=Locked  G=Email
$targetRow = 52

$dataRange = $ws.Range("A" + $targetRow + ":G" + $targetRow)
$dataRange.Borders.LineStyle = 1

$ws.Range("A" + $targetRow).Value = "Linking_AutoUser"
$ws.Range("B" + $targetRow).Value = "Password1"
$ws.Range("E" + $targetRow).Value = "Default user for Linking tests"
$ws.Range("F" + $targetRow).Value = "N"
$ws.Range("G" + $targetRow).Value = "linking.autouser@mailinator.com"

# Match the author's recorded view state: scrolled down with O31 selected.
$ws.Range("O31").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 1
